$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.951.63"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.233.25"
$ws.Range("E3").Value = "  +3.64%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'259.69"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").Value = "'82.89"
$ws.Range("E6").Value = "  +13.42%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "'44.41"
$ws.Range("E10").Value = "  +12.10%  "
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "'7.07"
$ws.Range("E12").Value = "  +5.06%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "2.567.88"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "2.245.24"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").Value = "43.859.43"
$ws.Range("E18").Value = "  +4.62%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "'71.17"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("D22").Value = "'2.38"
$ws.Range("E22").Value = "  +11.88%  "
$ws.Range("D23").Value = "'232.59"
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("D24").Value = "'9.30"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "'40.83"
$ws.Range("E27").Value = "  +10.74%  "
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'172.64"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").Value = "'0.0894"
$ws.Range("E32").Value = "  +12.19%  "
$ws.Range("D33").Value = "'20.68"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("E34").Value = "  +5.40%  "
$ws.Range("E35").Value = "  +9.30%  "
$ws.Range("D36").Value = "'0.0372"
$ws.Range("E36").Value = "  +13.30%  "
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("E38").Value = "  +7.02%  "
$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = "  +28.61%  "
$ws.Range("D40").Value = "'13.08"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").Value = "'63.65"
$ws.Range("E42").Value = "  +8.95%  "
$ws.Range("D43").Value = "'5.55"
$ws.Range("E43").Value = "  +8.29%  "
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "'104.34"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "'0.0988"
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("E48").Value = "  +30.22%  "
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("E51").Value = "  +4.06%  "
